$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / customer info
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number; force text so it isn't reinterpreted
# as a number, matching the original inline-string cell type.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 15.05.2025"

# Row 6
$ws.Range("B6").Value = "16.05."
$ws.Range("C6").Value = "17.05."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 57771522"
$ws.Range("E6").Value = "87,64-"

# Row 7
$ws.Range("B7").Value = "18.05."
$ws.Range("C7").Value = "19.05."
$ws.Range("D7").Value = "KARTENZ./18.05 ALDI SUED RO"
$ws.Range("E7").Value = "90,75-"

# Row 8
$ws.Range("B8").Value = "20.05."
$ws.Range("C8").Value = "21.05."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "24,57-"

# Row 9 (previously empty, now populated with a new transaction)
$ws.Range("B9").Value = "21.05."
$ws.Range("C9").Value = "22.05."
$ws.Range("D9").Value = "KARTENZ./21.05 LIDL RO"
$ws.Range("E9").Value = "24,21-"
$ws.Range("E9").HorizontalAlignment = -4152
$ws.Range("E9").VerticalAlignment = -4107
$ws.Range("E9").WrapText = $false

# Row 10 (previously empty, now populated with a new transaction)
$ws.Range("B10").Value = "22.05."
$ws.Range("C10").Value = "23.05."
$ws.Range("D10").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E10").Value = "56,77-"
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4107
$ws.Range("E10").WrapText = $false

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 27.05.2025"
$ws.Range("E12").Value = "283,94-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 01.06.2025"
